# Merge the split "<id>p107v_1</id>" runs back into a single run.
# The text is currently spread across three runs:
#   <id>   (Courier New, color 7f6000)
#   p107v_1 (default formatting)
#   </id>  (Courier New, color 7f6000)
# Find & Replace across that whole span collapses it into one run
# carrying the formatting of the first run in the matched range.
$d = $word.ActiveDocument
$d.Content.Find.Execute("<id>p107v_1</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p107v_1</id>", 2)
